# Add season "Wins" / "Losses" / "Ties" record columns (AD, AE, AF)
# to the existing player-stats table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Copy the formatting of the last existing header cell (AC1) onto the new
# header cells so they pick up the same bold/border/center style used by
# the rest of the header row, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-57) -------------------------------------------------
# Every player row gets the same team-season record.
$lastRow = 57
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 81   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 81   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
